$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AF4").Value = 0.694
$ws.Range("AF5").Value = 0.963
$ws.Range("AF6").Value = 0.8070000000000001
$ws.Range("AF7").Value = 0.894
$ws.Range("AF8").Value = 0.876
$ws.Range("AF9").Value = 0.741
$ws.Range("AF10").Value = 0.963
$ws.Range("AF11").Value = 0.963
$ws.Range("AF12").Value = 1.269
$ws.Range("AF13").Value = 1.667
